$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24, pushing the existing rows 24-30 down to 25-31.
$ws.Rows.Item(24).Insert()

# The constant columns (A,B,C,E,F,G,H,I,N,O,Q,R) are identical for every
# record in this block, so copy them from the row that used to be 24 and is
# now 25 into the newly inserted row 24.
$ws.Cells.Item(24, 1).Value2  = $ws.Cells.Item(25, 1).Value2   # A - Mercado ID
$ws.Cells.Item(24, 2).Value2  = $ws.Cells.Item(25, 2).Value2   # B - Mercado
$ws.Cells.Item(24, 3).Value2  = $ws.Cells.Item(25, 3).Value2   # C - Región
$ws.Cells.Item(24, 5).Value2  = $ws.Cells.Item(25, 5).Value2   # E - Codreg
$ws.Cells.Item(24, 6).Value2  = $ws.Cells.Item(25, 6).Value2   # F - Categoría ID
$ws.Cells.Item(24, 7).Value2  = $ws.Cells.Item(25, 7).Value2   # G - Categoría
$ws.Cells.Item(24, 8).Value2  = $ws.Cells.Item(25, 8).Value2   # H - Variedad
$ws.Cells.Item(24, 9).Value2  = $ws.Cells.Item(25, 9).Value2   # I - Calidad
$ws.Cells.Item(24, 14).Value2 = $ws.Cells.Item(25, 14).Value2  # N - Unidad de comercialización
$ws.Cells.Item(24, 15).Value2 = $ws.Cells.Item(25, 15).Value2  # O - Origen
$ws.Cells.Item(24, 17).Value2 = $ws.Cells.Item(25, 17).Value2  # Q - Kg o Unidades
$ws.Cells.Item(24, 18).Value2 = $ws.Cells.Item(25, 18).Value2  # R - Clasificación

# New weekly record values for row 24.
$ws.Cells.Item(24, 4).Value2  = 44582   # D - Fecha
$ws.Cells.Item(24, 10).Value2 = 520     # J - Volumen
$ws.Cells.Item(24, 11).Value2 = 15000   # K - Precio mínimo
$ws.Cells.Item(24, 12).Value2 = 16000   # L - Precio máximo
$ws.Cells.Item(24, 13).Value2 = 15500   # M - Precio promedio ponderado
$ws.Cells.Item(24, 16).Value2 = 1192    # P - Precio $/Kg

$ws.Range("D24").NumberFormat = $ws.Range("D25").NumberFormat
